$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"37.88856266666667"
$ws.Range("H2").Value = [double]"113.665688"
$ws.Range("I2").Value = [double]"0.9468476050819132"
$ws.Range("J2").Value = [double]"0.9540268599258594"
$ws.Range("M2").Value = [double]"162.7225033333333"
$ws.Range("N2").Value = [double]"488.16751"
$ws.Range("O2").Value = [double]"0.5231437953541009"
$ws.Range("P2").Value = [double]"0.5247717033381212"
$ws.Range("Q2").Value = [double]"6165.321764821875"
$ws.Range("R2").Value = [double]"55487.89588339688"
$ws.Range("S2").Value = [double]"0.495337449744493"
$ws.Range("T2").Value = [double]"0.5006463003136123"
$ws.Range("G3").Value = [double]"37.88856266666667"
$ws.Range("H3").Value = [double]"113.665688"
$ws.Range("I3").Value = [double]"0.9468476050819132"
$ws.Range("J3").Value = [double]"0.9540268599258594"
$ws.Range("O3").Value = [double]"0.0009322191998643353"
$ws.Range("P3").Value = [double]"0.0009351200601857102"
$ws.Range("Q3").Value = [double]"10.98633181459822"
$ws.Range("R3").Value = [double]"98.876986331384"
$ws.Range("S3").Value = [double]"0.0008826695168029234"
$ws.Range("T3").Value = [double]"0.0008921296546726537"
$ws.Range("G4").Value = [double]"37.88856266666667"
$ws.Range("H4").Value = [double]"113.665688"
$ws.Range("I4").Value = [double]"0.9468476050819132"
$ws.Range("J4").Value = [double]"0.9540268599258594"
$ws.Range("M4").Value = [double]"61.580654"
$ws.Range("N4").Value = [double]"184.741962"
$ws.Range("O4").Value = [double]"0.1979783766474813"
$ws.Range("P4").Value = [double]"0.1985944416431287"
$ws.Range("Q4").Value = [double]"2333.202468133317"
$ws.Range("R4").Value = [double]"20998.82221319986"
$ws.Range("S4").Value = [double]"0.1874553517866726"
$ws.Range("T4").Value = [double]"0.1894644315595234"
$ws.Range("G5").Value = [double]"37.88856266666667"
$ws.Range("H5").Value = [double]"113.665688"
$ws.Range("I5").Value = [double]"0.9468476050819132"
$ws.Range("J5").Value = [double]"0.9540268599258594"
$ws.Range("M5").Value = [double]"2.8947245"
$ws.Range("N5").Value = [double]"5.789449"
$ws.Range("O5").Value = [double]"0.009306378223129816"
$ws.Range("P5").Value = [double]"0.00622355841157717"
$ws.Range("Q5").Value = [double]"109.6769506209853"
$ws.Range("R5").Value = [double]"658.061703725912"
$ws.Range("S5").Value = [double]"0.008811721932556938"
$ws.Range("T5").Value = [double]"0.005937441888962137"
$ws.Range("G6").Value = [double]"37.88856266666667"
$ws.Range("H6").Value = [double]"113.665688"
$ws.Range("I6").Value = [double]"0.9468476050819132"
$ws.Range("J6").Value = [double]"0.9540268599258594"
$ws.Range("M6").Value = [double]"83.559527"
$ws.Range("N6").Value = [double]"250.678581"
$ws.Range("O6").Value = [double]"0.2686392305754237"
$ws.Range("P6").Value = [double]"0.2694751765469873"
$ws.Range("Q6").Value = [double]"3165.950375136526"
$ws.Range("R6").Value = [double]"28493.55337622873"
$ws.Range("S6").Value = [double]"0.2543604121013878"
$ws.Range("T6").Value = [double]"0.2570865565090888"
$ws.Range("I7").Value = [double]"0.02931771140176381"
$ws.Range("J7").Value = [double]"0.02954000622562442"
$ws.Range("M7").Value = [double]"162.7225033333333"
$ws.Range("N7").Value = [double]"488.16751"
$ws.Range("O7").Value = [double]"0.5231437953541009"
$ws.Range("P7").Value = [double]"0.5247717033381212"
$ws.Range("Q7").Value = [double]"190.8999116963744"
$ws.Range("R7").Value = [double]"1718.09920526737"
$ws.Range("S7").Value = [double]"0.01533737881381492"
$ws.Range("T7").Value = [double]"0.01550175938363963"
$ws.Range("I8").Value = [double]"0.02931771140176381"
$ws.Range("J8").Value = [double]"0.02954000622562442"
$ws.Range("O8").Value = [double]"0.0009322191998643353"
$ws.Range("P8").Value = [double]"0.0009351200601857102"
$ws.Range("S8").Value = [double]"2.733053346480576E-05"
$ws.Range("T8").Value = [double]"2.762345239959216E-05"
$ws.Range("I9").Value = [double]"0.02931771140176381"
$ws.Range("J9").Value = [double]"0.02954000622562442"
$ws.Range("M9").Value = [double]"61.580654"
$ws.Range("N9").Value = [double]"184.741962"
$ws.Range("O9").Value = [double]"0.1979783766474813"
$ws.Range("P9").Value = [double]"0.1985944416431287"
$ws.Range("Q9").Value = [double]"72.24410373483266"
$ws.Range("R9").Value = [double]"650.1969336134939"
$ws.Range("S9").Value = [double]"0.005804272910340552"
$ws.Range("T9").Value = [double]"0.005866481042512427"
$ws.Range("I10").Value = [double]"0.02931771140176381"
$ws.Range("J10").Value = [double]"0.02954000622562442"
$ws.Range("M10").Value = [double]"2.8947245"
$ws.Range("N10").Value = [double]"5.789449"
$ws.Range("O10").Value = [double]"0.009306378223129816"
$ws.Range("P10").Value = [double]"0.00622355841157717"
$ws.Range("Q10").Value = [double]"3.395981748777166"
$ws.Range("R10").Value = [double]"20.375890492663"
$ws.Range("S10").Value = [double]"0.0002728417109413794"
$ws.Range("T10").Value = [double]"0.0001838439542235268"
$ws.Range("I11").Value = [double]"0.02931771140176381"
$ws.Range("J11").Value = [double]"0.02954000622562442"
$ws.Range("M11").Value = [double]"83.559527"
$ws.Range("N11").Value = [double]"250.678581"
$ws.Range("O11").Value = [double]"0.2686392305754237"
$ws.Range("P11").Value = [double]"0.2694751765469873"
$ws.Range("Q11").Value = [double]"98.02888966754966"
$ws.Range("R11").Value = [double]"882.260007007947"
$ws.Range("S11").Value = [double]"0.007875887433202155"
$ws.Range("T11").Value = [double]"0.007960298392849244"
$ws.Range("G12").Value = [double]"0.05038133333333333"
$ws.Range("H12").Value = [double]"0.151144"
$ws.Range("I12").Value = [double]"0.00125904604054744"
$ws.Range("J12").Value = [double]"0.001268592468438093"
$ws.Range("M12").Value = [double]"162.7225033333333"
$ws.Range("N12").Value = [double]"488.16751"
$ws.Range("O12").Value = [double]"0.5231437953541009"
$ws.Range("P12").Value = [double]"0.5247717033381212"
$ws.Range("Q12").Value = [double]"8.198176681271111"
$ws.Range("R12").Value = [double]"73.78359013143999"
$ws.Range("S12").Value = [double]"0.0006586621241775413"
$ws.Range("T12").Value = [double]"0.0006657214305041697"
$ws.Range("G13").Value = [double]"0.05038133333333333"
$ws.Range("H13").Value = [double]"0.151144"
$ws.Range("I13").Value = [double]"0.00125904604054744"
$ws.Range("J13").Value = [double]"0.001268592468438093"
$ws.Range("O13").Value = [double]"0.0009322191998643353"
$ws.Range("P13").Value = [double]"0.0009351200601857102"
$ws.Range("Q13").Value = [double]"0.01460878973244444"
$ws.Range("R13").Value = [double]"0.131479107592"
$ws.Range("S13").Value = [double]"1.173706892511494E-06"
$ws.Range("T13").Value = [double]"1.186286265436968E-06"
$ws.Range("G14").Value = [double]"0.05038133333333333"
$ws.Range("H14").Value = [double]"0.151144"
$ws.Range("I14").Value = [double]"0.00125904604054744"
$ws.Range("J14").Value = [double]"0.001268592468438093"
$ws.Range("M14").Value = [double]"61.580654"
$ws.Range("N14").Value = [double]"184.741962"
$ws.Range("O14").Value = [double]"0.1979783766474813"
$ws.Range("P14").Value = [double]"0.1985944416431287"
$ws.Range("Q14").Value = [double]"3.102515456058667"
$ws.Range("R14").Value = [double]"27.922639104528"
$ws.Range("S14").Value = [double]"0.0002492638912320212"
$ws.Range("T14").Value = [double]"0.0002519354129421414"
$ws.Range("G15").Value = [double]"0.05038133333333333"
$ws.Range("H15").Value = [double]"0.151144"
$ws.Range("I15").Value = [double]"0.00125904604054744"
$ws.Range("J15").Value = [double]"0.001268592468438093"
$ws.Range("M15").Value = [double]"2.8947245"
$ws.Range("N15").Value = [double]"5.789449"
$ws.Range("O15").Value = [double]"0.009306378223129816"
$ws.Range("P15").Value = [double]"0.00622355841157717"
$ws.Range("Q15").Value = [double]"0.1458400799426667"
$ws.Range("R15").Value = [double]"0.8750404796560001"
$ws.Range("S15").Value = [double]"1.171715865366852E-05"
$ws.Range("T15").Value = [double]"7.895159327811339E-06"
$ws.Range("G16").Value = [double]"0.05038133333333333"
$ws.Range("H16").Value = [double]"0.151144"
$ws.Range("I16").Value = [double]"0.00125904604054744"
$ws.Range("J16").Value = [double]"0.001268592468438093"
$ws.Range("M16").Value = [double]"83.559527"
$ws.Range("N16").Value = [double]"250.678581"
$ws.Range("O16").Value = [double]"0.2686392305754237"
$ws.Range("P16").Value = [double]"0.2694751765469873"
$ws.Range("Q16").Value = [double]"4.209840382962667"
$ws.Range("R16").Value = [double]"37.888563446664"
$ws.Range("S16").Value = [double]"0.000338229159591698"
$ws.Range("T16").Value = [double]"0.0003418541793985335"
$ws.Range("G17").Value = [double]"0.903375"
$ws.Range("H17").Value = [double]"1.80675"
$ws.Range("I17").Value = [double]"0.02257563747577563"
$ws.Range("J17").Value = [double]"0.0151645413800781"
$ws.Range("M17").Value = [double]"162.7225033333333"
$ws.Range("N17").Value = [double]"488.16751"
$ws.Range("O17").Value = [double]"0.5231437953541009"
$ws.Range("P17").Value = [double]"0.5247717033381212"
$ws.Range("Q17").Value = [double]"146.99944144875"
$ws.Range("R17").Value = [double]"881.9966486925"
$ws.Range("S17").Value = [double]"0.01181030467161554"
$ws.Range("T17").Value = [double]"0.007957922210365008"
$ws.Range("G18").Value = [double]"0.903375"
$ws.Range("H18").Value = [double]"1.80675"
$ws.Range("I18").Value = [double]"0.02257563747577563"
$ws.Range("J18").Value = [double]"0.0151645413800781"
$ws.Range("O18").Value = [double]"0.0009322191998643353"
$ws.Range("P18").Value = [double]"0.0009351200601857102"
$ws.Range("Q18").Value = [double]"0.261946529625"
$ws.Range("R18").Value = [double]"1.57167917775"
$ws.Range("S18").Value = [double]"2.104544270409486E-05"
$ws.Range("T18").Value = [double]"1.418066684802732E-05"
$ws.Range("G19").Value = [double]"0.903375"
$ws.Range("H19").Value = [double]"1.80675"
$ws.Range("I19").Value = [double]"0.02257563747577563"
$ws.Range("J19").Value = [double]"0.0151645413800781"
$ws.Range("M19").Value = [double]"61.580654"
$ws.Range("N19").Value = [double]"184.741962"
$ws.Range("O19").Value = [double]"0.1979783766474813"
$ws.Range("P19").Value = [double]"0.1985944416431287"
$ws.Range("Q19").Value = [double]"55.63042330725001"
$ws.Range("R19").Value = [double]"333.7825398435"
$ws.Range("S19").Value = [double]"0.004469488059236102"
$ws.Range("T19").Value = [double]"0.003011593628150731"
$ws.Range("G20").Value = [double]"0.903375"
$ws.Range("H20").Value = [double]"1.80675"
$ws.Range("I20").Value = [double]"0.02257563747577563"
$ws.Range("J20").Value = [double]"0.0151645413800781"
$ws.Range("M20").Value = [double]"2.8947245"
$ws.Range("N20").Value = [double]"5.789449"
$ws.Range("O20").Value = [double]"0.009306378223129816"
$ws.Range("P20").Value = [double]"0.00622355841157717"
$ws.Range("Q20").Value = [double]"2.6150217451875"
$ws.Range("R20").Value = [double]"10.46008698075"
$ws.Range("S20").Value = [double]"0.0002100974209778317"
$ws.Range("T20").Value = [double]"9.437740906369513E-05"
$ws.Range("G21").Value = [double]"0.903375"
$ws.Range("H21").Value = [double]"1.80675"
$ws.Range("I21").Value = [double]"0.02257563747577563"
$ws.Range("J21").Value = [double]"0.0151645413800781"
$ws.Range("M21").Value = [double]"83.559527"
$ws.Range("N21").Value = [double]"250.678581"
$ws.Range("O21").Value = [double]"0.2686392305754237"
$ws.Range("P21").Value = [double]"0.2694751765469873"
$ws.Range("Q21").Value = [double]"75.485587703625"
$ws.Range("R21").Value = [double]"452.91352622175"
$ws.Range("S21").Value = [double]"0.006064701881242065"
$ws.Range("T21").Value = [double]"0.004086467465650639"
